$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-20 01:50:01"

for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
